# The deck's theme (ppt/theme/theme1.xml, linked from the slide master)
# currently holds the "Integral" color scheme. The target edit swaps its
# colors for the stock "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink), matching what ppt/theme/theme2.xml already contains.
#
# PowerPoint's RGB() helper isn't available in this host, so pack the
# R/G/B triplet into the same little-endian integer VBA's RGB() produces.
function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2
$colors.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2
$colors.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1
$colors.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2
$colors.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3
$colors.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4
$colors.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5
$colors.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6
$colors.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink
$colors.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink
